# Insert a new weekly data row at row 208, shifting existing rows 208-314
# down to 209-315 (matches the Hortaliza / Ciboulette weekly update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 208:314 down by inserting a new blank row at 208.
$ws.Rows("208:208").Insert()

# Populate the newly inserted row 208 with the new weekly record.
$ws.Cells.Item(208, 1).Value = 6
$ws.Cells.Item(208, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(208, 3).Value = "Metropolitana"
$ws.Cells.Item(208, 4).Value = 44510
$ws.Cells.Item(208, 5).Value = 13
$ws.Cells.Item(208, 6).Value = 100112039
$ws.Cells.Item(208, 7).Value = "Ciboulette"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 820
$ws.Cells.Item(208, 11).Value = 700
$ws.Cells.Item(208, 12).Value = 800
$ws.Cells.Item(208, 13).Value = 746
$ws.Cells.Item(208, 14).Value = "$/docena de atados"
$ws.Cells.Item(208, 15).Value = "Región Metropolitana"
$ws.Cells.Item(208, 16).Value = 249
$ws.Cells.Item(208, 17).Value = 3
$ws.Cells.Item(208, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(208, 4).NumberFormat = $ws.Cells.Item(209, 4).NumberFormat
